$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5620
$ws1.Range("F4").Value = 642
$ws1.Range("F6").Value = 834
$ws1.Range("F7").Value = 54
$ws1.Range("F9").Value = 6
$ws1.Range("F10").Value = 6
$ws1.Range("F11").Value = 20

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 20

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5620
$ws4.Range("F4").Value = 642
$ws4.Range("F6").Value = 834
$ws4.Range("F7").Value = 54
$ws4.Range("F10").Value = 6
$ws4.Range("F11").Value = 6
$ws4.Range("F12").Value = 20
$ws4.Range("F13").Value = 20
